# New Microsite scripts support to Beta server
# Adds new registration-history rows to the "AMSIN" sheet and the "AMS" sheet,
# and corrects an existing timestamp on "AMS" row 15.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write one data row (columns A..G) into a worksheet, matching the
# existing column layout:
#   A = Run Date (text, yyyy-mm-dd)
#   B = Run Time (datetime serial, custom date/time format)
#   C = Sprint Name (text)
#   D = Total Cases (number)
#   E = Pass Cases (number)
#   F = Fail Cases (number)
#   G = Time Taken (number)
# ---------------------------------------------------------------------------

function Write-HistoryRow {
    param(
        $ws,
        $row,
        $styleRow,
        $applyStyle,
        $runDate,
        $runTime,
        $sprintName,
        $totalCases,
        $passCases,
        $failCases,
        $timeTaken
    )

    # Column A - Run Date, stored as literal text (not auto-converted to a date serial)
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $runDate
    if ($applyStyle) {
        $ws.Cells.Item($row, 1).Style = $ws.Cells.Item($styleRow, 1).Style
    } else {
        $ws.Cells.Item($row, 1).ClearFormats()
    }

    # Column B - Run Time, numeric datetime serial using the same number format as the reference row
    $ws.Cells.Item($row, 2).Style = $ws.Cells.Item($styleRow, 2).Style
    $ws.Cells.Item($row, 2).NumberFormat = $ws.Cells.Item($styleRow, 2).NumberFormat
    $ws.Cells.Item($row, 2).Value = $runTime

    # Column C - Sprint Name (text)
    $ws.Cells.Item($row, 3).Value = $sprintName
    if ($applyStyle) {
        $ws.Cells.Item($row, 3).Style = $ws.Cells.Item($styleRow, 3).Style
    }

    # Column D - Total Cases
    $ws.Cells.Item($row, 4).Value = $totalCases
    if ($applyStyle) {
        $ws.Cells.Item($row, 4).Style = $ws.Cells.Item($styleRow, 4).Style
    }

    # Column E - Pass Cases
    $ws.Cells.Item($row, 5).Value = $passCases
    if ($applyStyle) {
        $ws.Cells.Item($row, 5).Style = $ws.Cells.Item($styleRow, 5).Style
    }

    # Column F - Fail Cases
    $ws.Cells.Item($row, 6).Value = $failCases
    if ($applyStyle) {
        $ws.Cells.Item($row, 6).Style = $ws.Cells.Item($styleRow, 6).Style
    }

    # Column G - Time Taken
    $ws.Cells.Item($row, 7).Value = $timeTaken
    if ($applyStyle) {
        $ws.Cells.Item($row, 7).Style = $ws.Cells.Item($styleRow, 7).Style
    }
}

# ---------------------------------------------------------------------------
# Sheet "AMSIN": append rows 23-26 (new Beta/microsite OCR registration runs)
# ---------------------------------------------------------------------------

$wsAmsin = $wb.Worksheets.Item("AMSIN")

Write-HistoryRow $wsAmsin 23 22 $true "2022-09-15" 44819.71062266204 "ocrecs166"   42 42 0 1.44
Write-HistoryRow $wsAmsin 24 22 $true "2022-09-16" 44820.64606804398 "fstcocr167"  42 42 0 1.36
Write-HistoryRow $wsAmsin 25 22 $true "2022-09-19" 44823.60064732639 "scndocr167"  42 42 0 2.03
Write-HistoryRow $wsAmsin 26 22 $true "2022-09-20" 44824.38858300926 "finalocr167" 42 42 0 1.77

# ---------------------------------------------------------------------------
# Sheet "AMS": fix the run time recorded on row 15, give row 15 the same
# explicit styling as the rows above it, and append the new Beta row 16.
# ---------------------------------------------------------------------------

$wsAms = $wb.Worksheets.Item("AMS")

# Row 15 already holds 2022-09-08 / ocr166 data; only the Run Time value
# changes, while every cell in the row picks up explicit styling (matching
# row 14's style) instead of the workbook default.
$wsAms.Cells.Item(15, 1).Style = $wsAms.Cells.Item(14, 1).Style
$wsAms.Cells.Item(15, 2).Value = 44812.54125284722
$wsAms.Cells.Item(15, 3).Style = $wsAms.Cells.Item(14, 3).Style
$wsAms.Cells.Item(15, 4).Style = $wsAms.Cells.Item(14, 4).Style
$wsAms.Cells.Item(15, 5).Style = $wsAms.Cells.Item(14, 5).Style
$wsAms.Cells.Item(15, 6).Style = $wsAms.Cells.Item(14, 6).Style
$wsAms.Cells.Item(15, 7).Style = $wsAms.Cells.Item(14, 7).Style

# New row 16 keeps the workbook's default (unstyled) formatting, same as row
# 15 had before this edit - only the Run Time column carries a number format.
Write-HistoryRow $wsAms 16 15 $false "2022-09-20" 44824.73127158471 "betaocr167" 42 42 0 1.13
